$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Formula = "=12-3"
$ws.Range("A7").Select()
